$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.707324028015137
$ws.Range("B1").Value = 1.891999006271362
$ws.Range("C1").Value = 2.26327657699585
$ws.Range("D1").Value = 3.560611724853516
$ws.Range("E1").Value = 1.977408170700073
